# Refactoring 9/29/24 @ 13:02
# Remove the "confirmationCode" column (F) from the Password sheet.
# This shifts the accessToken/refreshToken columns (G,H) one column to the
# left (into F,G) and removes the now-empty trailing column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Password")

# Delete the confirmationCode column's cells (F1:F2) and shift everything
# in F:H to the left by one column.
$ws.Range("F1:F2").Delete(-4159)  # xlShiftToLeft
